# feat: add 2022-Q1 data
#
# The workbook has sheets: 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计
# We need to insert a new "2022-Q1" sheet (fund holding detail, same layout
# as the existing quarter sheets) right before "总计", and refresh "总计"
# (the quarterly roll-up) with a new first data row for 2022-Q1 plus the
# renumbered index column.
#
# Strategy: rename the existing "总计" sheet to "2022-Q1" (so it keeps the
# sheetId/position continuity and the header/index cell styling already on
# it), then add a brand new sheet right after it and rename that one "总计".
# This reproduces the exact sheet order / sheetId allocation seen in the
# target workbook.

$wb = $excel.ActiveWorkbook

$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# ---------------------------------------------------------------------
# 1. "2022-Q1" sheet: fund holding detail table (B:H, rows 1-13)
# ---------------------------------------------------------------------

# Extend the existing bold/bordered header style (already on B1:D1 from the
# old "总计" sheet) across the new columns, and the index-column style
# (already on A2:A6) down across the new rows, by copying formats from the
# cells that already carry it - this reuses the existing style entries
# instead of minting new ones.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$q1.Range("A7:A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$fundData = @(
    @("000689", "前海开源新经济灵活配置混合", "142.62", "94.48", "3.84", "5.4766", 7),
    @("010490", "鹏华高质量增长混合A", "13.31", "93.61", "8.83", "1.1753", 1),
    @("009023", "鹏华稳健回报混合", "3.52", "93.91", "9.21", "0.3242", 2),
    @("005535", "泰信竞争优选灵活配置混合", "10.07", "90.21", "3.03", "0.3051", 10),
    @("008381", "前海开源新兴产业混合", "6.18", "93.63", "4.84", "0.2991", 10),
    @("290006", "泰信蓝筹精选混合", "7.62", "88.52", "2.99", "0.2278", 10),
    @("014036", "博时成长回报混合A", "5.50", "68.14", "2.69", "0.1480", 7),
    @("014037", "博时成长回报混合C", "1.00", "68.14", "2.69", "0.0269", 7),
    @("010491", "鹏华高质量增长混合C", "0.28", "93.61", "8.83", "0.0247", 1),
    @("002495", "前海开源量化优选灵活配置混合A", "0.49", "93.24", "3.16", "0.0155", 3),
    @("002496", "前海开源量化优选灵活配置混合C", "0.26", "93.24", "3.16", "0.0082", 3),
    @("970083", "东海证券海盈6个月持有期混合", "0.14", "20.43", "1.38", "0.0019", 7)
)

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $i + 2
    $row = $fundData[$i]
    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1.Cells.Item($r, 6).Value = "'" + $row[4]
    $q1.Cells.Item($r, 7).Value = "'" + $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2. "总计" sheet: quarterly roll-up table (B:D, rows 1-7)
# ---------------------------------------------------------------------

# Bring over the header / index-column styling from the "2022-Q1" sheet
# (same reuse-not-mint approach as above, just cross-sheet).
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @("2022-Q1", 12, 8.029999999999999),
    @("2021-Q4", 12, 9),
    @("2021-Q3", 1, 3.97),
    @("2021-Q2", 3, 2.78),
    @("2021-Q1", 7, 0.76),
    @("2020-Q4", 28, 15.1)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $r = $i + 2
    $row = $totalData[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}
